{"js": "// The document contains one table of two-digit-division problems laid out\n// as 5 \"content\" rows (each holding 5 answer cells) interleaved with blank\n// spacer rows. The commit swaps each answer cell's text for a new problem,\n// in strict reading order (row major, left-to-right, top-to-bottom).\n//\n// Because a couple of the \"before\" strings repeat (e.g. \"81\u00f74=20, 1\"\n// appears twice but maps to two different replacements), we must match by\n// position (row/column index), not by a global text search-and-replace.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row indices (within the 20-row table) that actually hold answer text;\n// the rows in between are empty spacer rows.\nconst contentRows = [0, 4, 8, 12, 16];\n\n// Replacement text for each cell, in row-major order matching the diff.\nconst replacements = [\n  [\"70\u00f74=17, 2\", \"11\u00f77=1, 4\"],\n  [\"82\u00f72=41, 0\", \"71\u00f74=17, 3\"],\n  [\"72\u00f72=36, 0\", \"95\u00f79=10, 5\"],\n  [\"79\u00f74=19, 3\", \"62\u00f73=20, 2\"],\n  [\"27\u00f76=4, 3\", \"26\u00f75=5, 1\"],\n  [\"81\u00f74=20, 1\", \"88\u00f75=17, 3\"],\n  [\"45\u00f74=11, 1\", \"18\u00f73=6, 0\"],\n  [\"64\u00f79=7, 1\", \"67\u00f76=11, 1\"],\n  [\"27\u00f74=6, 3\", \"36\u00f79=4, 0\"],\n  [\"43\u00f73=14, 1\", \"20\u00f78=2, 4\"],\n  [\"44\u00f77=6, 2\", \"68\u00f72=34, 0\"],\n  [\"68\u00f75=13, 3\", \"74\u00f75=14, 4\"],\n  [\"28\u00f75=5, 3\", \"78\u00f74=19, 2\"],\n  [\"79\u00f78=9, 7\", \"99\u00f72=49, 1\"],\n  [\"94\u00f73=31, 1\", \"62\u00f73=20, 2\"],\n  [\"14\u00f77=2, 0\", \"83\u00f72=41, 1\"],\n  [\"69\u00f72=34, 1\", \"67\u00f79=7, 4\"],\n  [\"57\u00f74=14, 1\", \"90\u00f74=22, 2\"],\n  [\"84\u00f73=28, 0\", \"58\u00f72=29, 0\"],\n  [\"57\u00f72=28, 1\", \"52\u00f74=13, 0\"],\n  [\"77\u00f78=9, 5\", \"52\u00f77=7, 3\"],\n  [\"45\u00f75=9, 0\", \"74\u00f74=18, 2\"],\n  [\"45\u00f76=7, 3\", \"18\u00f72=9, 0\"],\n  [\"46\u00f75=9, 1\", \"19\u00f73=6, 1\"],\n  [\"81\u00f74=20, 1\", \"61\u00f78=7, 5\"],\n];\n\nlet k = 0;\nfor (const rowIndex of contentRows) {\n  for (let col = 0; col < 5; col++) {\n    const cell = table.getCellOrNullObject(rowIndex, col);\n    cell.load(\"value\");\n    await context.sync();\n\n    const [expectedOld, newValue] = replacements[k];\n    if (!cell.isNullObject && cell.value === expectedOld) {\n      cell.value = newValue;\n    }\n    k++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains one table of two-digit-division problems laid out\n# as 5 \"content\" rows (each holding 5 answer cells) interleaved with blank\n# spacer rows. The commit swaps each answer cell's text for a new problem,\n# in strict reading order (row major, left-to-right, top-to-bottom).\n#\n# Because a couple of the \"before\" strings repeat (e.g. \"81\u00f74=20, 1\"\n# appears twice but maps to two different replacements), we must match by\n# position (row/column index), not by a global text search-and-replace.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1-based row indices (within the 20-row table) that hold answer text; the\n# rows in between are empty spacer rows.\n$contentRows = @(1, 5, 9, 13, 17)\n\n# Replacement text for each cell, in row-major order matching the diff.\n$replacements = @(\n    @(\"70\u00f74=17, 2\", \"11\u00f77=1, 4\"),\n    @(\"82\u00f72=41, 0\", \"71\u00f74=17, 3\"),\n    @(\"72\u00f72=36, 0\", \"95\u00f79=10, 5\"),\n    @(\"79\u00f74=19, 3\", \"62\u00f73=20, 2\"),\n    @(\"27\u00f76=4, 3\", \"26\u00f75=5, 1\"),\n    @(\"81\u00f74=20, 1\", \"88\u00f75=17, 3\"),\n    @(\"45\u00f74=11, 1\", \"18\u00f73=6, 0\"),\n    @(\"64\u00f79=7, 1\", \"67\u00f76=11, 1\"),\n    @(\"27\u00f74=6, 3\", \"36\u00f79=4, 0\"),\n    @(\"43\u00f73=14, 1\", \"20\u00f78=2, 4\"),\n    @(\"44\u00f77=6, 2\", \"68\u00f72=34, 0\"),\n    @(\"68\u00f75=13, 3\", \"74\u00f75=14, 4\"),\n    @(\"28\u00f75=5, 3\", \"78\u00f74=19, 2\"),\n    @(\"79\u00f78=9, 7\", \"99\u00f72=49, 1\"),\n    @(\"94\u00f73=31, 1\", \"62\u00f73=20, 2\"),\n    @(\"14\u00f77=2, 0\", \"83\u00f72=41, 1\"),\n    @(\"69\u00f72=34, 1\", \"67\u00f79=7, 4\"),\n    @(\"57\u00f74=14, 1\", \"90\u00f74=22, 2\"),\n    @(\"84\u00f73=28, 0\", \"58\u00f72=29, 0\"),\n    @(\"57\u00f72=28, 1\", \"52\u00f74=13, 0\"),\n    @(\"77\u00f78=9, 5\", \"52\u00f77=7, 3\"),\n    @(\"45\u00f75=9, 0\", \"74\u00f74=18, 2\"),\n    @(\"45\u00f76=7, 3\", \"18\u00f72=9, 0\"),\n    @(\"46\u00f75=9, 1\", \"19\u00f73=6, 1\"),\n    @(\"81\u00f74=20, 1\", \"61\u00f78=7, 5\")\n)\n\n$k = 0\nforeach ($row in $contentRows) {\n    for ($col = 1; $col -le 5; $col++) {\n        $cell = $t.Cell($row, $col)\n        $pair = $replacements[$k]\n        $expectedOld = $pair[0]\n        $newValue = $pair[1]\n\n        # Cell.Range.Text includes trailing cell-mark characters (CR + BEL);\n        # strip them before comparing against the expected source text.\n        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($current -eq $expectedOld) {\n            $cell.Range.Text = $newValue\n        }\n        $k++\n    }\n}\n"}
